# Corrected excel sheets for application fix issues
#
# Summary of the real edits applied to this workbook:
#  - Summary sheet: a new (blank, unstyled) cell G2 appears next to the
#    existing data block, and the sheet's selection moves to rows 7-15.
#  - Repayment schedule sheet: a new column O (mirroring column N's blank /
#    zero pattern) is added for rows 2-15, and the sheet's selection moves
#    to row 16.
#  - Transactions sheet: three loan/transaction IDs are corrected
#    (54->90, 53->89, 52->88) and the sheet's selection moves to D3.
#  - The active tab moves from "Repayment schedule" to "Transactions".

$wb = $excel.ActiveWorkbook

$wsSummary    = $wb.Worksheets.Item("Summary")
$wsRepayment  = $wb.Worksheets.Item("Repayment schedule")
$wsTxn        = $wb.Worksheets.Item("Transactions")

# ---------------------------------------------------------------------
# Summary sheet: extend the used range with a blank G2 cell (copy the
# format of an already-bare cell so it stays styleless/valueless, just
# like the real edit), then move the selection down to rows 7:15.
# ---------------------------------------------------------------------
$null = $wsTxn.Range("K4").Copy()
$null = $wsSummary.Range("G2").PasteSpecial(-4122)   # xlPasteFormats

$null = $wsSummary.Activate()
$null = $wsSummary.Range("A7:A15").EntireRow.Select()

# ---------------------------------------------------------------------
# Repayment schedule sheet: add column O (same look as column N) across
# rows 2-15 -- blank on rows 2 & 4, zero everywhere else -- then move the
# selection down to row 16.
# ---------------------------------------------------------------------
$null = $wsRepayment.Range("N2:N15").Copy()
$null = $wsRepayment.Range("O2:O15").PasteSpecial(-4122)   # xlPasteFormats

$wsRepayment.Range("O3").Value = 0
$wsRepayment.Range("O5:O15").Value = 0

$null = $wsRepayment.Activate()
$null = $wsRepayment.Range("A16:A16").EntireRow.Select()

# ---------------------------------------------------------------------
# Transactions sheet: correct the ID values and move the selection to D3.
# This also becomes the active sheet/tab, matching the workbook's
# activeTab moving from the Repayment schedule tab to this one.
# ---------------------------------------------------------------------
$wsTxn.Range("A2").Value = 90
$wsTxn.Range("A3").Value = 89
$wsTxn.Range("A4").Value = 88

$null = $wsTxn.Activate()
$null = $wsTxn.Range("D3").Select()
